# "hacked excels with 23/4 added" - append the 23-April-2019 data row
# (row 68) to the bottom of the table, mirroring the date-cell styling
# used by the rest of column A.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 68

# Copy row 67's formatting into row 68 first so the new date cell (A68)
# picks up the same style (built-in date number format) as every other
# cell in column A, then overwrite the values.
$ws.Range("A67").Copy($ws.Range("A68"))

$ws.Cells.Item($row, 1).Value = 43578
$ws.Cells.Item($row, 2).Value = 4021
$ws.Cells.Item($row, 3).Value = 643
$ws.Cells.Item($row, 4).Value = 214
$ws.Cells.Item($row, 5).Value = 1012
$ws.Cells.Item($row, 6).Value = 3024
$ws.Cells.Item($row, 7).Value = 2641

# Leave the view looking like the new row was just entered.
$ws.Range("G70").Select()
